$d = $word.ActiveDocument

function Fix-FigureCaption($oldPrefix, $newText) {
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $capPara = $d.Paragraphs.Item($i)
        $capText = $capPara.Range.Text
        if ($capText -like "$oldPrefix*") {
            # The paragraph immediately before the caption holds the inline figure.
            $imgPara = $capPara.Previous()
            $imgPara.Style = "Captioned Figure"

            $capRange = $capPara.Range
            $capRange.End = $capRange.End - 1
            $capRange.Text = ""
            $capRange.Font.Bold = $false
            $capRange.Text = $newText
            $capPara.Style = "Image Caption"
            return
        }
    }
}

Fix-FigureCaption "Figure S1:" "Figure 1 Scatter plot of XXX. Each point indicates XXXX."
Fix-FigureCaption "Figure S2:" "Figure 2 Boxplot of XXXX."

Write-Output "done"
